$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9, pushing existing rows 9.. down by one.
$ws.Rows.Item(9).Insert()

# Fill the newly inserted row 9 with the new record's data.
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44616
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 100114007
$ws.Cells.Item(9, 7).Value = "Jengibre"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 22
$ws.Cells.Item(9, 11).Value = 25000
$ws.Cells.Item(9, 12).Value = 26000
$ws.Cells.Item(9, 13).Value = 25545
$ws.Cells.Item(9, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 1965
$ws.Cells.Item(9, 17).Value = 13
$ws.Cells.Item(9, 18).Value = "Hortaliza"
